$d = $word.ActiveDocument

$replacements = @(
    @{old = "766÷5=153, 1"; new = "420÷2=210, 0"},
    @{old = "958÷5=191, 3"; new = "687÷8=85, 7"},
    @{old = "600÷9=66, 6"; new = "240÷6=40, 0"},
    @{old = "335÷6=55, 5"; new = "623÷6=103, 5"},
    @{old = "595÷3=198, 1"; new = "391÷6=65, 1"},
    @{old = "978÷8=122, 2"; new = "253÷2=126, 1"},
    @{old = "582÷2=291, 0"; new = "719÷9=79, 8"},
    @{old = "129÷4=32, 1"; new = "317÷3=105, 2"},
    @{old = "526÷7=75, 1"; new = "530÷8=66, 2"},
    @{old = "780÷7=111, 3"; new = "859÷9=95, 4"},
    @{old = "441÷3=147, 0"; new = "803÷2=401, 1"},
    @{old = "717÷2=358, 1"; new = "969÷8=121, 1"},
    @{old = "980÷6=163, 2"; new = "157÷3=52, 1"},
    @{old = "302÷7=43, 1"; new = "658÷6=109, 4"},
    @{old = "282÷2=141, 0"; new = "880÷8=110, 0"},
    @{old = "984÷6=164, 0"; new = "271÷8=33, 7"},
    @{old = "820÷5=164, 0"; new = "265÷8=33, 1"},
    @{old = "759÷8=94, 7"; new = "682÷9=75, 7"},
    @{old = "157÷9=17, 4"; new = "155÷3=51, 2"},
    @{old = "875÷3=291, 2"; new = "672÷5=134, 2"},
    @{old = "188÷2=94, 0"; new = "317÷6=52, 5"},
    @{old = "563÷9=62, 5"; new = "544÷8=68, 0"},
    @{old = "489÷6=81, 3"; new = "352÷8=44, 0"},
    @{old = "710÷5=142, 0"; new = "215÷2=107, 1"},
    @{old = "634÷9=70, 4"; new = "408÷8=51, 0"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
